$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.816.45'
$ws.Range('E2').Value = '  -1.00%  '
$ws.Range('D3').Value = '1.941.88'
$ws.Range('E3').Value = '  -0.92%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '242.08'
$ws.Range('E5').Value = '  -2.06%  '
$ws.Range('E6').Value = '  -0.10%  '
$ws.Range('D7').Value = '0.4887'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '0.2954'
$ws.Range('E8').Value = '  -0.78%  '
$ws.Range('D9').Value = '0.06893'
$ws.Range('E9').Value = '  +0.75%  '
$ws.Range('D10').Value = '19.44'
$ws.Range('E10').Value = '  +1.54%  '
$ws.Range('D11').Value = '106.30'
$ws.Range('E11').Value = '  -0.24%  '
$ws.Range('D12').Value = '1.957.01'
$ws.Range('E12').Value = '  -0.05%  '
$ws.Range('D13').Value = '0.07724'
$ws.Range('E13').Value = '  -0.42%  '
$ws.Range('D14').Value = '5.345'
$ws.Range('E14').Value = '  -1.52%  '
$ws.Range('D15').Value = '0.6993'
$ws.Range('E15').Value = '  -1.91%  '
$ws.Range('D16').Value = '277.07'
$ws.Range('E16').Value = '  -2.92%  '
$ws.Range('D17').Value = '30.813.60'
$ws.Range('E17').Value = '  -1.03%  '
$ws.Range('D18').Value = '0.000007727'
$ws.Range('E18').Value = '  -0.67%  '
$ws.Range('D19').Value = '13.11'
$ws.Range('E19').Value = '  -1.06%  '
$ws.Range('B20').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C20').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D20').Value = '2.197.27'
$ws.Range('E20').Value = '  +0.40%  '
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').Value = '0.9998'
$ws.Range('E21').Value = '  -0.12%  '
$ws.Range('D22').Value = '5.445'
$ws.Range('E22').Value = '  -2.59%  '
$ws.Range('D23').Value = '0.9992'
$ws.Range('E23').Value = '  -0.17%  '
$ws.Range('D24').Value = '6.514'
$ws.Range('E24').Value = '  -1.34%  '
$ws.Range('D25').Value = '9.716'
$ws.Range('E25').Value = '  -2.84%  '
$ws.Range('D26').Value = '167.80'
$ws.Range('E26').Value = '  -0.58%  '
$ws.Range('D27').Value = '19.65'
$ws.Range('E27').Value = '  -1.95%  '
$ws.Range('D28').Value = '2.155'
$ws.Range('E28').Value = '  -2.11%  '
$ws.Range('D29').Value = '0.1045'
$ws.Range('E29').Value = '  -1.62%  '
$ws.Range('D30').Value = '1.387'
$ws.Range('E30').Value = '  -4.02%  '
$ws.Range('D31').Value = '1.551'
$ws.Range('E31').Value = '  -2.82%  '
$ws.Range('D32').Value = '4.550'
$ws.Range('E32').Value = '  -5.54%  '
$ws.Range('E33').Value = '  -3.22%  '
$ws.Range('D34').Value = '0.04852'
$ws.Range('E34').Value = '  -3.53%  '
$ws.Range('D35').Value = '0.7512'
$ws.Range('E35').Value = '  -2.49%  '
$ws.Range('D36').Value = '1.158'
$ws.Range('E36').Value = '  -0.77%  '
$ws.Range('D37').Value = '0.9995'
$ws.Range('E37').Value = '  -0.10%  '
$ws.Range('E38').Value = '  -0.37%  '
$ws.Range('D39').Value = '0.01988'
$ws.Range('E40').Value = '  -1.94%  '
$ws.Range('D41').Value = '78.67'
$ws.Range('D42').Value = '6.450'
$ws.Range('E42').Value = '  +0.46%  '
$ws.Range('D43').Value = '2.096'
$ws.Range('E43').Value = '  -1.74%  '
$ws.Range('D44').Value = '0.9126'
$ws.Range('E44').Value = '  +2.95%  '
$ws.Range('D45').Value = '108.80'
$ws.Range('E45').Value = '  -0.83%  '
$ws.Range('D46').Value = '0.4400'
$ws.Range('E46').Value = '  -1.59%  '
$ws.Range('E47').Value = '  -0.33%  '
$ws.Range('D48').Value = '7.721'
$ws.Range('E48').Value = '  +2.87%  '
$ws.Range('D49').Value = '983.73'
$ws.Range('E49').Value = '  -1.47%  '
$ws.Range('D50').Value = '0.1244'
$ws.Range('E50').Value = '  -2.09%  '
$ws.Range('D51').Value = '9.285'
$ws.Range('E51').Value = '  -1.33%  '
